$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value2 = 12
$ws.Cells.Item(2, 2).Value2 = 'h$_{q}$'
$ws.Cells.Item(2, 3).Value2 = 0.00251885861143194
$ws.Cells.Item(3, 1).Value2 = 34
$ws.Cells.Item(3, 2).Value2 = '$F_{q}$'
$ws.Cells.Item(3, 3).Value2 = 0.001589514856865449
$ws.Cells.Item(4, 1).Value2 = 60
$ws.Cells.Item(4, 2).Value2 = '$(F_{p}^{\text{SCF}})_{3}$'
$ws.Cells.Item(4, 3).Value2 = 0.0004958128213381861
$ws.Cells.Item(5, 1).Value2 = 19
$ws.Cells.Item(5, 2).Value2 = '(h$_{rs}$)$_{1}$'
$ws.Cells.Item(5, 3).Value2 = 0.0004606031230109178
$ws.Cells.Item(6, 1).Value2 = 43
$ws.Cells.Item(6, 2).Value2 = '$\eta_{s}$'
$ws.Cells.Item(6, 3).Value2 = 0.0004599159301036932
$ws.Cells.Item(7, 1).Value2 = 32
$ws.Cells.Item(7, 2).Value2 = '$F_{q}^{\text{SCF}}$'
$ws.Cells.Item(7, 3).Value2 = 0.0004402905619992127
$ws.Cells.Item(8, 1).Value2 = 40
$ws.Cells.Item(8, 2).Value2 = '$F_{s}^{\text{SCF}}$'
$ws.Cells.Item(8, 3).Value2 = 0.0003695023415982953
$ws.Cells.Item(9, 1).Value2 = 73
$ws.Cells.Item(9, 2).Value2 = '$\langle ss \vert ss \rangle$'
$ws.Cells.Item(9, 3).Value2 = 0.0003507350350169923
$ws.Cells.Item(10, 1).Value2 = 13
$ws.Cells.Item(10, 2).Value2 = 'h$_{qs}$'
$ws.Cells.Item(10, 3).Value2 = 0.0003465515291880913
$ws.Cells.Item(11, 1).Value2 = 59
$ws.Cells.Item(11, 2).Value2 = '$(\eta_{r})_{2}$'
$ws.Cells.Item(11, 3).Value2 = 0.0003158557496052796
$ws.Cells.Item(12, 1).Value2 = 52
$ws.Cells.Item(12, 2).Value2 = '$(F_{p}^{\text{SCF}})_{2}$'
$ws.Cells.Item(12, 3).Value2 = 0.0002854065069971326
$ws.Cells.Item(13, 1).Value2 = 71
$ws.Cells.Item(13, 2).Value2 = '$\langle qq \vert qq \rangle$'
$ws.Cells.Item(13, 3).Value2 = 0.0002633924884706824
$ws.Cells.Item(14, 1).Value2 = 41
$ws.Cells.Item(14, 2).Value2 = '$\omega_{s}$'
$ws.Cells.Item(14, 3).Value2 = 0.0002450684975358547
$ws.Cells.Item(15, 1).Value2 = 16
$ws.Cells.Item(15, 2).Value2 = '(h$_{r}$)$_{2}$'
$ws.Cells.Item(15, 3).Value2 = 0.0001657487616176912
$ws.Cells.Item(16, 1).Value2 = 99
$ws.Cells.Item(16, 2).Value2 = '$(\langle pq \vert qp \rangle)_{3}$'
$ws.Cells.Item(16, 3).Value2 = 0.0001372035474603177
$ws.Cells.Item(17, 1).Value2 = 7
$ws.Cells.Item(17, 2).Value2 = '(h$_{pq}$)$_{3}$'
$ws.Cells.Item(17, 3).Value2 = 0.0001257939397919092
$ws.Cells.Item(18, 1).Value2 = 101
$ws.Cells.Item(18, 2).Value2 = '$(\langle rs \vert sr \rangle)_{3}$'
$ws.Cells.Item(18, 3).Value2 = 0.0001168103239308631
$ws.Cells.Item(19, 1).Value2 = 1
$ws.Cells.Item(19, 2).Value2 = '(h$_{p}$)$_{1}$'
$ws.Cells.Item(19, 3).Value2 = 0.000114511546037373
$ws.Cells.Item(20, 1).Value2 = 2
$ws.Cells.Item(20, 2).Value2 = '(h$_{p}$)$_{2}$'
$ws.Cells.Item(20, 3).Value2 = 0.0001057437360674874
$ws.Cells.Item(21, 1).Value2 = 64
$ws.Cells.Item(21, 2).Value2 = '$(F_{r}^{\text{SCF}})_{3}$'
$ws.Cells.Item(21, 3).Value2 = 0.000100490705045045
$ws.Cells.Item(22, 1).Value2 = 22
$ws.Cells.Item(22, 2).Value2 = 'h$_{s}$'
$ws.Cells.Item(22, 3).Value2 = 0.00009732365029823314
$ws.Cells.Item(23, 1).Value2 = 98
$ws.Cells.Item(23, 2).Value2 = '$(\langle pq \vert pq \rangle)_{3}$'
$ws.Cells.Item(23, 3).Value2 = 0.00009704512267078167
$ws.Cells.Item(24, 1).Value2 = 97
$ws.Cells.Item(24, 2).Value2 = '$(\langle rr \vert rr \rangle)_{3}$'
$ws.Cells.Item(24, 3).Value2 = 0.00009565283581139257
$ws.Cells.Item(25, 1).Value2 = 17
$ws.Cells.Item(25, 2).Value2 = '(h$_{r}$)$_{3}$'
$ws.Cells.Item(25, 3).Value2 = 0.00009538472743242263
$ws.Cells.Item(26, 1).Value2 = 42
$ws.Cells.Item(26, 2).Value2 = '$F_{s}$'
$ws.Cells.Item(26, 3).Value2 = 0.00009537455866418268
$ws.Cells.Item(27, 1).Value2 = 0
$ws.Cells.Item(27, 2).Value2 = '(h$_{p}$)$_{0}$'
$ws.Cells.Item(27, 3).Value2 = 0.00008969583424564923
$ws.Cells.Item(28, 1).Value2 = 23
$ws.Cells.Item(28, 2).Value2 = 'type_0'
$ws.Cells.Item(28, 3).Value2 = 0.00008167506535324721
$ws.Cells.Item(29, 1).Value2 = 11
$ws.Cells.Item(29, 2).Value2 = '(h$_{pr}$)$_{3}$'
$ws.Cells.Item(29, 3).Value2 = 0.00008105753693794668
$ws.Cells.Item(30, 1).Value2 = 3
$ws.Cells.Item(30, 2).Value2 = '(h$_{p}$)$_{3}$'
$ws.Cells.Item(30, 3).Value2 = 0.00007897983315358066
$ws.Cells.Item(31, 1).Value2 = 20
$ws.Cells.Item(31, 2).Value2 = '(h$_{rs}$)$_{2}$'
$ws.Cells.Item(31, 3).Value2 = 0.00007886413641980995
$ws.Cells.Item(32, 1).Value2 = 66
$ws.Cells.Item(32, 2).Value2 = '$(F_{r})_{3}$'
$ws.Cells.Item(32, 3).Value2 = 0.0000658118714152279
$ws.Cells.Item(33, 1).Value2 = 75
$ws.Cells.Item(33, 2).Value2 = '$(\langle pq \vert qp \rangle)_{0}$'
$ws.Cells.Item(33, 3).Value2 = 0.00006430601206920176
$ws.Cells.Item(34, 1).Value2 = 18
$ws.Cells.Item(34, 2).Value2 = '(h$_{rs}$)$_{0}$'
$ws.Cells.Item(34, 3).Value2 = 0.00006427183390137861
$ws.Cells.Item(35, 1).Value2 = 100
$ws.Cells.Item(35, 2).Value2 = '$(\langle rs\vert rs \rangle)_{3}$'
$ws.Cells.Item(35, 3).Value2 = 0.0000608770675349818
$ws.Cells.Item(36, 1).Value2 = 80
$ws.Cells.Item(36, 2).Value2 = '$(\langle pp \vert pp \rangle)_{1}$'
$ws.Cells.Item(36, 3).Value2 = 0.00005734212437046565
$ws.Cells.Item(37, 1).Value2 = 15
$ws.Cells.Item(37, 2).Value2 = '(h$_{r}$)$_{1}$'
$ws.Cells.Item(37, 3).Value2 = 0.00005654982150563682
$ws.Cells.Item(38, 1).Value2 = 91
$ws.Cells.Item(38, 2).Value2 = '$(\langle pq \vert qp \rangle)_{2}$'
$ws.Cells.Item(38, 3).Value2 = 0.00005607008682399571
$ws.Cells.Item(39, 1).Value2 = 5
$ws.Cells.Item(39, 2).Value2 = '(h$_{pq}$)$_{1}$'
$ws.Cells.Item(39, 3).Value2 = 0.00005509515003640848
$ws.Cells.Item(40, 1).Value2 = 25
$ws.Cells.Item(40, 2).Value2 = 'type_2'
$ws.Cells.Item(40, 3).Value2 = 0.00005267716249619125
$ws.Cells.Item(41, 1).Value2 = 56
$ws.Cells.Item(41, 2).Value2 = '$(F_{r}^{\text{SCF}})_{2}$'
$ws.Cells.Item(41, 3).Value2 = 0.00005156309668670854
$ws.Cells.Item(42, 1).Value2 = 89
$ws.Cells.Item(42, 2).Value2 = '$(\langle rr \vert rr \rangle)_{2}$'
$ws.Cells.Item(42, 3).Value2 = 0.00005044488406731033
$ws.Cells.Item(43, 1).Value2 = 26
$ws.Cells.Item(43, 2).Value2 = 'type_3'
$ws.Cells.Item(43, 3).Value2 = 0.00004852964924115941
$ws.Cells.Item(44, 1).Value2 = 92
$ws.Cells.Item(44, 2).Value2 = '$(\langle rs\vert rs \rangle)_{2}$'
$ws.Cells.Item(44, 3).Value2 = 0.00004769573412697785
$ws.Cells.Item(45, 1).Value2 = 10
$ws.Cells.Item(45, 2).Value2 = '(h$_{pr}$)$_{2}$'
$ws.Cells.Item(45, 3).Value2 = 0.00004701906598103767
$ws.Cells.Item(46, 1).Value2 = 67
$ws.Cells.Item(46, 2).Value2 = '$(\eta_{r})_{3}$'
$ws.Cells.Item(46, 3).Value2 = 0.00004549863048493938
$ws.Cells.Item(47, 1).Value2 = 21
$ws.Cells.Item(47, 2).Value2 = '(h$_{rs}$)$_{3}$'
$ws.Cells.Item(47, 3).Value2 = 0.00004294800809414763
$ws.Cells.Item(48, 1).Value2 = 28
$ws.Cells.Item(48, 2).Value2 = '$(F_{p}^{\text{SCF}})_{0}$'
$ws.Cells.Item(48, 3).Value2 = 0.00004280454055400444
$ws.Cells.Item(49, 1).Value2 = 82
$ws.Cells.Item(49, 2).Value2 = '$(\langle pq \vert pq \rangle)_{1}$'
$ws.Cells.Item(49, 3).Value2 = 0.00004034940450682245
$ws.Cells.Item(50, 1).Value2 = 96
$ws.Cells.Item(50, 2).Value2 = '$(\langle pp \vert pp \rangle)_{3}$'
$ws.Cells.Item(50, 3).Value2 = 0.00003686082926544592
$ws.Cells.Item(51, 1).Value2 = 36
$ws.Cells.Item(51, 2).Value2 = '$(F_{r}^{\text{SCF}})_{0}$'
$ws.Cells.Item(51, 3).Value2 = 0.00003321065450982063
$ws.Cells.Item(52, 1).Value2 = 83
$ws.Cells.Item(52, 2).Value2 = '$(\langle pq \vert qp \rangle)_{1}$'
$ws.Cells.Item(52, 3).Value2 = 0.00003240910889744105
$ws.Cells.Item(53, 1).Value2 = 88
$ws.Cells.Item(53, 2).Value2 = '$(\langle pp \vert pp \rangle)_{2}$'
$ws.Cells.Item(53, 3).Value2 = 0.00003170832436738017
$ws.Cells.Item(54, 1).Value2 = 85
$ws.Cells.Item(54, 2).Value2 = '$(\langle rs \vert sr \rangle)_{1}$'
$ws.Cells.Item(54, 3).Value2 = 0.00003082676075309285
$ws.Cells.Item(55, 1).Value2 = 4
$ws.Cells.Item(55, 2).Value2 = '(h$_{pq}$)$_{0}$'
$ws.Cells.Item(55, 3).Value2 = 0.00002985623103002895
$ws.Cells.Item(56, 1).Value2 = 6
$ws.Cells.Item(56, 2).Value2 = '(h$_{pq}$)$_{2}$'
$ws.Cells.Item(56, 3).Value2 = 0.00002949753579171774
$ws.Cells.Item(57, 1).Value2 = 74
$ws.Cells.Item(57, 2).Value2 = '$(\langle pq \vert pq \rangle)_{0}$'
$ws.Cells.Item(57, 3).Value2 = 0.00002851108909777619
$ws.Cells.Item(58, 1).Value2 = 44
$ws.Cells.Item(58, 2).Value2 = '$(F_{p}^{\text{SCF}})_{1}$'
$ws.Cells.Item(58, 3).Value2 = 0.00002659829157757205
$ws.Cells.Item(59, 1).Value2 = 46
$ws.Cells.Item(59, 2).Value2 = '$(F_{p})_{1}$'
$ws.Cells.Item(59, 3).Value2 = 0.00002654548416042313
$ws.Cells.Item(60, 1).Value2 = 33
$ws.Cells.Item(60, 2).Value2 = '$\omega_{q}$'
$ws.Cells.Item(60, 3).Value2 = 0.00002603161923351117
$ws.Cells.Item(61, 1).Value2 = 84
$ws.Cells.Item(61, 2).Value2 = '$(\langle rs\vert rs \rangle)_{1}$'
$ws.Cells.Item(61, 3).Value2 = 0.0000242572823892681
$ws.Cells.Item(62, 1).Value2 = 48
$ws.Cells.Item(62, 2).Value2 = '$(F_{r}^{\text{SCF}})_{1}$'
$ws.Cells.Item(62, 3).Value2 = 0.00002084932871852011
$ws.Cells.Item(63, 1).Value2 = 93
$ws.Cells.Item(63, 2).Value2 = '$(\langle rs \vert sr \rangle)_{2}$'
$ws.Cells.Item(63, 3).Value2 = 0.00002032336059824532
$ws.Cells.Item(64, 1).Value2 = 54
$ws.Cells.Item(64, 2).Value2 = '$(F_{p})_{2}$'
$ws.Cells.Item(64, 3).Value2 = 0.00002007449324700505
$ws.Cells.Item(65, 1).Value2 = 86
$ws.Cells.Item(65, 2).Value2 = '$(\langle pq \vert rs \rangle)_{2}$'
$ws.Cells.Item(65, 3).Value2 = 0.00001888256016577377
$ws.Cells.Item(66, 1).Value2 = 76
$ws.Cells.Item(66, 2).Value2 = '$(\langle rs\vert rs \rangle)_{0}$'
$ws.Cells.Item(66, 3).Value2 = 0.00001882737731857092
$ws.Cells.Item(67, 1).Value2 = 62
$ws.Cells.Item(67, 2).Value2 = '$(F_{p})_{3}$'
$ws.Cells.Item(67, 3).Value2 = 0.00001729672510891181
$ws.Cells.Item(68, 1).Value2 = 77
$ws.Cells.Item(68, 2).Value2 = '$(\langle rs \vert sr \rangle)_{0}$'
$ws.Cells.Item(68, 3).Value2 = 0.00001563981510390768
$ws.Cells.Item(69, 1).Value2 = 72
$ws.Cells.Item(69, 2).Value2 = '$(\langle rr \vert rr \rangle)_{0}$'
$ws.Cells.Item(69, 3).Value2 = 0.0000155342879402321
$ws.Cells.Item(70, 1).Value2 = 94
$ws.Cells.Item(70, 2).Value2 = '$(\langle pq \vert rs \rangle)_{3}$'
$ws.Cells.Item(70, 3).Value2 = 0.00001543432891697375
$ws.Cells.Item(71, 1).Value2 = 81
$ws.Cells.Item(71, 2).Value2 = '$(\langle rr \vert rr \rangle)_{1}$'
$ws.Cells.Item(71, 3).Value2 = 0.00001446734971568919
$ws.Cells.Item(72, 1).Value2 = 78
$ws.Cells.Item(72, 2).Value2 = '$(\langle pq \vert rs \rangle)_{1}$'
$ws.Cells.Item(72, 3).Value2 = 0.00001383353897681905
$ws.Cells.Item(73, 1).Value2 = 14
$ws.Cells.Item(73, 2).Value2 = '(h$_{r}$)$_{0}$'
$ws.Cells.Item(73, 3).Value2 = 0.00001236922424315693
$ws.Cells.Item(74, 1).Value2 = 9
$ws.Cells.Item(74, 2).Value2 = '(h$_{pr}$)$_{1}$'
$ws.Cells.Item(74, 3).Value2 = 0.00001177016027270752
$ws.Cells.Item(75, 1).Value2 = 8
$ws.Cells.Item(75, 2).Value2 = '(h$_{pr}$)$_{0}$'
$ws.Cells.Item(75, 3).Value2 = 0.000009453233038684323
$ws.Cells.Item(76, 1).Value2 = 90
$ws.Cells.Item(76, 2).Value2 = '$(\langle pq \vert pq \rangle)_{2}$'
$ws.Cells.Item(76, 3).Value2 = 0.000008837189700135434
$ws.Cells.Item(77, 1).Value2 = 70
$ws.Cells.Item(77, 2).Value2 = '$(\langle pp \vert pp \rangle)_{0}$'
$ws.Cells.Item(77, 3).Value2 = 0.000007901779935853004
$ws.Cells.Item(78, 1).Value2 = 58
$ws.Cells.Item(78, 2).Value2 = '$(F_{r})_{2}$'
$ws.Cells.Item(78, 3).Value2 = 0.000007398604825193421
$ws.Cells.Item(79, 1).Value2 = 24
$ws.Cells.Item(79, 2).Value2 = 'type_1'
$ws.Cells.Item(79, 3).Value2 = 0.000007217731296794699
$ws.Cells.Item(80, 1).Value2 = 51
$ws.Cells.Item(80, 2).Value2 = '$(\eta_{r})_{1}$'
$ws.Cells.Item(80, 3).Value2 = 0.000005449617962603447
$ws.Cells.Item(81, 1).Value2 = 50
$ws.Cells.Item(81, 2).Value2 = '$(F_{r})_{1}$'
$ws.Cells.Item(81, 3).Value2 = 0.000005271179289827433
$ws.Cells.Item(82, 1).Value2 = 55
$ws.Cells.Item(82, 2).Value2 = '$(\eta_{p})_{2}$'
$ws.Cells.Item(82, 3).Value2 = 0.000004685141095172145
$ws.Cells.Item(83, 1).Value2 = 63
$ws.Cells.Item(83, 2).Value2 = '$(\eta_{p})_{3}$'
$ws.Cells.Item(83, 3).Value2 = 0.000004608585150642206
$ws.Cells.Item(84, 1).Value2 = 68
$ws.Cells.Item(84, 2).Value2 = '$(\langle pq \vert rs \rangle)_{0}$'
$ws.Cells.Item(84, 3).Value2 = 0.000004027767037629824
$ws.Cells.Item(85, 1).Value2 = 30
$ws.Cells.Item(85, 2).Value2 = '$(F_{p})_{0}$'
$ws.Cells.Item(85, 3).Value2 = 0.000003394123468384157
$ws.Cells.Item(86, 1).Value2 = 39
$ws.Cells.Item(86, 2).Value2 = '$(\eta_{r})_{0}$'
$ws.Cells.Item(86, 3).Value2 = 0.000003191839366313202
$ws.Cells.Item(87, 1).Value2 = 35
$ws.Cells.Item(87, 2).Value2 = '$\eta_{q}$'
$ws.Cells.Item(87, 3).Value2 = 0.00000311016522193979
$ws.Cells.Item(88, 1).Value2 = 31
$ws.Cells.Item(88, 2).Value2 = '$(\eta_{p})_{0}$'
$ws.Cells.Item(88, 3).Value2 = 0.000002738589303588591
$ws.Cells.Item(89, 1).Value2 = 47
$ws.Cells.Item(89, 2).Value2 = '$(\eta_{p})_{1}$'
$ws.Cells.Item(89, 3).Value2 = 0.000002310746494203963
$ws.Cells.Item(90, 1).Value2 = 38
$ws.Cells.Item(90, 2).Value2 = '$(F_{r})_{0}$'
$ws.Cells.Item(90, 3).Value2 = 0.000001967942964748721
$ws.Cells.Item(91, 1).Value2 = 95
$ws.Cells.Item(91, 2).Value2 = '$(\langle pq \vert sr \rangle)_{3}$'
$ws.Cells.Item(91, 3).Value2 = 0.00000002360379197396587
$ws.Cells.Item(92, 1).Value2 = 49
$ws.Cells.Item(92, 2).Value2 = '$(\omega_{r})_{1}$'
$ws.Cells.Item(92, 3).Value2 = 0.00000002264392748862727
$ws.Cells.Item(93, 1).Value2 = 65
$ws.Cells.Item(93, 2).Value2 = '$(\omega_{r})_{3}$'
$ws.Cells.Item(93, 3).Value2 = 0.00000002221166636860165
$ws.Cells.Item(94, 1).Value2 = 69
$ws.Cells.Item(94, 2).Value2 = '$(\langle pq \vert sr \rangle)_{0}$'
$ws.Cells.Item(94, 3).Value2 = 0.00000002127593006271969
$ws.Cells.Item(95, 1).Value2 = 79
$ws.Cells.Item(95, 2).Value2 = '$(\langle pq \vert sr \rangle)_{1}$'
$ws.Cells.Item(95, 3).Value2 = 0.00000001891066643230072
$ws.Cells.Item(96, 1).Value2 = 87
$ws.Cells.Item(96, 2).Value2 = '$(\langle pq \vert sr \rangle)_{2}$'
$ws.Cells.Item(96, 3).Value2 = 0.00000001789212365922618
$ws.Cells.Item(97, 1).Value2 = 53
$ws.Cells.Item(97, 2).Value2 = '$(\omega_{p})_{2}$'
$ws.Cells.Item(97, 3).Value2 = 0.00000001537510843060651
$ws.Cells.Item(98, 1).Value2 = 27
$ws.Cells.Item(98, 2).Value2 = '$\mathbf{b}$'
$ws.Cells.Item(98, 3).Value2 = 0.00000001314790846253896
$ws.Cells.Item(99, 1).Value2 = 37
$ws.Cells.Item(99, 2).Value2 = '$(\omega_{r})_{0}$'
$ws.Cells.Item(99, 3).Value2 = 0.00000001210670940526213
$ws.Cells.Item(100, 1).Value2 = 57
$ws.Cells.Item(100, 2).Value2 = '$(\omega_{r})_{2}$'
$ws.Cells.Item(100, 3).Value2 = 0.000000008872060532500205
$ws.Cells.Item(101, 1).Value2 = 29
$ws.Cells.Item(101, 2).Value2 = '$(\omega_{p})_{0}$'
$ws.Cells.Item(101, 3).Value2 = 0.000000008634345780532474
$ws.Cells.Item(102, 1).Value2 = 61
$ws.Cells.Item(102, 2).Value2 = '$(\omega_{p})_{3}$'
$ws.Cells.Item(102, 3).Value2 = 0.000000006569235958214897
$ws.Cells.Item(103, 1).Value2 = 45
$ws.Cells.Item(103, 2).Value2 = '$(\omega_{p})_{1}$'
$ws.Cells.Item(103, 3).Value2 = 0.000000003270562620510968

Write-Output "Updated rows"
